# Auto update Excel log
# Appends newly-logged sensor rows to the ALERTS, PIR, Humidity, and
# Temperature sheets, matching new sensor readings captured on 2026-02-06.

function Add-LogRow {
    param($ws, $row, $date, $timestamp, $hour, $location, $value, $status)

    # Several columns hold text that looks numeric to Excel's auto-detection
    # (dates like "2026-02-06", percentages like "70.6%", etc.). Left alone,
    # Excel would silently convert such text into a date serial / percentage
    # number and attach a number-format style to the cell. The source log
    # stores everything as plain text with no cell styling, so force the
    # whole row range to Text format first, write the literal values, then
    # clear the formatting back to the workbook default (no explicit style)
    # to match the rest of the sheet.
    $rng = $ws.Range("A" + $row + ":F" + $row)
    $rng.NumberFormat = "@"

    $ws.Range("A$row").Value = $date
    $ws.Range("B$row").Value = $timestamp
    $ws.Range("C$row").Value = $hour
    $ws.Range("D$row").Value = $location
    $ws.Range("E$row").Value = $value
    $ws.Range("F$row").Value = $status

    $rng.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALERTS sheet: add rows 11-12
# ---------------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $wsAlerts 11 "2026-02-06" "09:56:09" "09:00" "Bathroom" "MODERATE" "MODERATE ALERT: Bathroom occupied, no motion > 40s."
Add-LogRow $wsAlerts 12 "2026-02-06" "09:56:30" "09:00" "Bathroom" "CRITICAL" "CRITICAL ALERT: Bathroom occupied, no motion > 60s."

# ---------------------------------------------------------------------------
# PIR sheet: add rows 192-204
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRow $wsPIR 192 "2026-02-06" "09:55:57" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 193 "2026-02-06" "09:55:59" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 194 "2026-02-06" "09:56:04" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 195 "2026-02-06" "09:56:10" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 196 "2026-02-06" "09:56:14" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 197 "2026-02-06" "09:56:20" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 198 "2026-02-06" "09:56:25" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 199 "2026-02-06" "09:56:30" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 200 "2026-02-06" "09:56:35" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 201 "2026-02-06" "09:56:40" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 202 "2026-02-06" "09:56:45" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 203 "2026-02-06" "09:56:50" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPIR 204 "2026-02-06" "09:56:55" "09:00" "Bathroom" "No Motion" "Inactive"

# ---------------------------------------------------------------------------
# Humidity sheet: add rows 110-121
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRow $wsHumidity 110 "2026-02-06" "09:55:58" "09:00" "Bathroom" "70.6%" "Active"
Add-LogRow $wsHumidity 111 "2026-02-06" "09:56:02" "09:00" "Bathroom" "70.6%" "Active"
Add-LogRow $wsHumidity 112 "2026-02-06" "09:56:07" "09:00" "Bathroom" "70.6%" "Active"
Add-LogRow $wsHumidity 113 "2026-02-06" "09:56:12" "09:00" "Bathroom" "70.5%" "Active"
Add-LogRow $wsHumidity 114 "2026-02-06" "09:56:18" "09:00" "Bathroom" "70.4%" "Active"
Add-LogRow $wsHumidity 115 "2026-02-06" "09:56:23" "09:00" "Bathroom" "70.3%" "Active"
Add-LogRow $wsHumidity 116 "2026-02-06" "09:56:27" "09:00" "Bathroom" "70.2%" "Active"
Add-LogRow $wsHumidity 117 "2026-02-06" "09:56:32" "09:00" "Bathroom" "70.1%" "Active"
Add-LogRow $wsHumidity 118 "2026-02-06" "09:56:38" "09:00" "Bathroom" "70.2%" "Active"
Add-LogRow $wsHumidity 119 "2026-02-06" "09:56:43" "09:00" "Bathroom" "70.2%" "Active"
Add-LogRow $wsHumidity 120 "2026-02-06" "09:56:48" "09:00" "Bathroom" "70.2%" "Active"
Add-LogRow $wsHumidity 121 "2026-02-06" "09:56:53" "09:00" "Bathroom" "69.2%" "Active"

# ---------------------------------------------------------------------------
# Temperature sheet: add rows 110-121
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRow $wsTemperature 110 "2026-02-06" "09:55:58" "09:00" "Bathroom" "27.8C" "Active"
Add-LogRow $wsTemperature 111 "2026-02-06" "09:56:03" "09:00" "Bathroom" "27.8C" "Active"
Add-LogRow $wsTemperature 112 "2026-02-06" "09:56:08" "09:00" "Bathroom" "27.8C" "Active"
Add-LogRow $wsTemperature 113 "2026-02-06" "09:56:13" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 114 "2026-02-06" "09:56:18" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 115 "2026-02-06" "09:56:23" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 116 "2026-02-06" "09:56:28" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 117 "2026-02-06" "09:56:33" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 118 "2026-02-06" "09:56:38" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 119 "2026-02-06" "09:56:43" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 120 "2026-02-06" "09:56:48" "09:00" "Bathroom" "27.9C" "Active"
Add-LogRow $wsTemperature 121 "2026-02-06" "09:56:53" "09:00" "Bathroom" "27.9C" "Active"
